# This workbook's weekly data refresh reshuffles the existing data rows
# (rows 2-26, columns A:T) into a new row order. Row 15 stays in place.
# Build a map of target row -> source row (both refer to the ORIGINAL
# layout before any writes happen), snapshot all source rows first, then
# write them out in the new order so overlapping reads/writes don't
# clobber each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 26
$firstCol = 1   # A
$lastCol = 20   # T

# target row -> source row (original positions)
$rowMap = @{
    2  = 26
    3  = 16
    4  = 17
    5  = 2
    6  = 3
    7  = 12
    8  = 19
    9  = 13
    10 = 25
    11 = 24
    12 = 6
    13 = 7
    14 = 8
    15 = 15
    16 = 22
    17 = 23
    18 = 10
    19 = 11
    20 = 14
    21 = 20
    22 = 9
    23 = 4
    24 = 18
    25 = 21
    26 = 5
}

# Snapshot every source row's values (A:T) before any write occurs.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowRange = $ws.Range($ws.Cells.Item($r, $firstCol), $ws.Cells.Item($r, $lastCol))
    $snapshot[$r] = $rowRange.Value()
}

# Now write each target row using the snapshotted source row values.
for ($targetRow = $firstRow; $targetRow -le $lastRow; $targetRow++) {
    $sourceRow = $rowMap[$targetRow]
    $values = $snapshot[$sourceRow]
    $destRange = $ws.Range($ws.Cells.Item($targetRow, $firstCol), $ws.Cells.Item($targetRow, $lastCol))
    $destRange.Value = $values
}
